$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of (id, speaker_variant) pairs for rows 2..20.
# Column D (is_prefered) is cleared for every data row, and a new
# row (20) is appended at the bottom.
$data = @(
    @("#wolfs-aert", "Wolfs aert"),
    @("#verkeerde-meyning", "Verkeerde meyning"),
    @("#verkeerde-meyninghe", "Verkeerde meyninghe"),
    @("#goed-aert", "Goed aert"),
    @("#verkeerde-meyninge", "Verkeerde meyninge"),
    @("#verkeerde-meyning", "verkeerde meyning"),
    @("#sondich-bedrijf", "Sondich bedrijf"),
    @("#ootmoet", "Ootmoet"),
    @("#menich-slechtaert", "Menich slechtaert"),
    @("#goet-onderwijs", "Goet onderwijs"),
    @("#'thert-vol-vreesen", "'Thert vol vreesen"),
    @("#magdalena", "Magdalena"),
    @("#vvolfs-aert", "VVolfs aert"),
    @("#verkeerde-meyningh", "Verkeerde meyningh"),
    @("#thert-vol-vreesen", "Thert vol vreesen"),
    @("#onverduldicheyt", "Onverduldicheyt"),
    @("#schriftueren-troost", "Schriftueren troost"),
    @("#verkeerde-meyningh", "verkeerde meyningh"),
    @("#g", "g")
)

$row = 2
foreach ($pair in $data) {
    $idCell = $ws.Cells.Item($row, 2)
    $textCell = $ws.Cells.Item($row, 3)

    $idVal = $pair[0]
    $textVal = $pair[1]

    # A value beginning with a literal apostrophe needs to be doubled so
    # Excel doesn't treat it as a text-qualifier prefix and strip it; the
    # style is then reset so no stray quote-prefix formatting is left
    # behind on the cell.
    if ($idVal.StartsWith("'")) {
        $idCell.Value = "'" + $idVal
        $idCell.Style = "Normal"
    } else {
        $idCell.Value = $idVal
    }

    if ($textVal.StartsWith("'")) {
        $textCell.Value = "'" + $textVal
        $textCell.Style = "Normal"
    } else {
        $textCell.Value = $textVal
    }

    # Clear is_prefered, but keep the (now-empty) cell present.
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = ""
    $dCell.Style = "Normal"

    $row = $row + 1
}

# Fill column A (URL) for the newly added row 20.
$ws.Range("A20").Value = "https://www.dbnl.org/tekst/heyn003pest01_01"

# Row 20 is brand new, so E20:H20 need to be created as empty (but
# present) cells too, matching the rest of the table's layout.
for ($c = 5; $c -le 8; $c++) {
    $cell = $ws.Cells.Item(20, $c)
    $cell.Value = ""
    $cell.Style = "Normal"
}
